$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = "'69.355.57"
$ws.Range('E2').Value = '  -0.12%  '

# Row 3
$ws.Range('D3').Value = "'3.670.55"
$ws.Range('E3').Value = '  -0.49%  '

# Row 4
$ws.Range('D4').Value = "'1.00"
$ws.Range('E4').Value = '  +0.20%  '

# Row 5
$ws.Range('D5').Value = "'644.47"
$ws.Range('E5').Value = '  -5.41%  '

# Row 6
$ws.Range('D6').Value = "'158.59"
$ws.Range('E6').Value = '  -1.07%  '

# Row 7
$ws.Range('D7').Value = "'1.00"
$ws.Range('E7').Value = '  +0.02%  '

# Row 8
$ws.Range('E8').Value = '  +0.36%  '

# Row 9
$ws.Range('D9').Value = "'0.144"
$ws.Range('E9').Value = '  -1.26%  '

# Row 10
$ws.Range('D10').Value = "'7.07"
$ws.Range('E10').Value = '  -1.48%  '

# Row 11
$ws.Range('D11').Value = "'0.440"
$ws.Range('E11').Value = '  +0.21%  '

# Row 12
$ws.Range('D12').Value = "'0.0000230"
$ws.Range('E12').Value = '  -1.14%  '

# Row 13
$ws.Range('D13').Value = "'4.297.63"
$ws.Range('E13').Value = '  -0.24%  '

# Row 14
$ws.Range('D14').Value = "'32.42"
$ws.Range('E14').Value = '  +0.25%  '

# Row 15
$ws.Range('D15').Value = "'3.669.97"
$ws.Range('E15').Value = '  -0.69%  '

# Row 16
$ws.Range('D16').Value = "'69.440.22"
$ws.Range('E16').Value = '  +0.10%  '

# Row 17
$ws.Range('E17').Value = '  -0.14%  '

# Row 18
$ws.Range('D18').Value = "'15.86"
$ws.Range('E18').Value = '  -0.91%  '

# Row 19
$ws.Range('D19').Value = "'6.44"
$ws.Range('E19').Value = '  -0.42%  '

# Row 20
$ws.Range('D20').Value = "'467.07"
$ws.Range('E20').Value = '  -1.08%  '

# Row 21
$ws.Range('D21').Value = "'10.00"
$ws.Range('E21').Value = '  +1.49%  '

# Row 22
$ws.Range('D22').Value = "'0.644"
$ws.Range('E22').Value = '  -0.90%  '

# Row 23
$ws.Range('D23').Value = "'79.36"
$ws.Range('E23').Value = '  -0.99%  '

# Row 24
$ws.Range('D24').Value = "'3.823.64"
$ws.Range('E24').Value = '  -0.27%  '

# Row 25
$ws.Range('E25').Value = '  +0.08%  '

# Row 26
$ws.Range('D26').Value = "'0.0000124"
$ws.Range('E26').Value = '  -0.15%  '

# Row 27
$ws.Range('D27').Value = "'10.68"
$ws.Range('E27').Value = '  -2.03%  '

# Row 28
$ws.Range('D28').Value = "'8.96"
$ws.Range('E28').Value = '  -1.61%  '

# Row 29
$ws.Range('D29').Value = "'2.62"
$ws.Range('E29').Value = '  -3.27%  '

# Row 30
$ws.Range('D30').Value = "'1.68"
$ws.Range('E30').Value = '  -2.83%  '

# Row 31
$ws.Range('D31').Value = "'1.99"
$ws.Range('E31').Value = '  -0.16%  '

# Row 32
$ws.Range('D32').Value = "'1.00"
$ws.Range('E32').Value = '  -0.20%  '

# Row 33
$ws.Range('D33').Value = "'26.89"
$ws.Range('E33').Value = '  -0.22%  '

# Row 34
$ws.Range('D34').Value = "'6.40"
$ws.Range('E34').Value = '  -2.19%  '

# Row 35
$ws.Range('B35').Value = 'RenzoRestakedETH'
$ws.Range('C35').Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range('D35').Value = "'3.669.72"
$ws.Range('E35').Value = '  -0.15%  '

# Row 36
$ws.Range('B36').Value = 'Kaspa'
$ws.Range('C36').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D36').Value = "'0.162"
$ws.Range('E36').Value = '  -0.33%  '

# Row 37
$ws.Range('D37').Value = "'8.39"
$ws.Range('E37').Value = '  +0.00%  '

# Row 39
$ws.Range('B39').Value = 'Monero'
$ws.Range('C39').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D39').Value = "'178.85"
$ws.Range('E39').Value = '  +5.53%  '

# Row 40
$ws.Range('B40').Value = 'FirstDigitalUSD'
$ws.Range('C40').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D40').Value = "'1.00"
$ws.Range('E40').Value = '  +0.20%  '

# Row 41
$ws.Range('B41').Value = 'Filecoin'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D41').Value = "'5.83"
$ws.Range('E41').Value = '  -6.70%  '

# Row 42
$ws.Range('D42').Value = "'2.21"
$ws.Range('E42').Value = '  -2.67%  '

# Row 43
$ws.Range('D43').Value = "'0.0889"
$ws.Range('E43').Value = '  -1.88%  '

# Row 44
$ws.Range('D44').Value = "'0.923"
$ws.Range('E44').Value = '  -1.84%  '

# Row 45
$ws.Range('D45').Value = "'47.14"
$ws.Range('E45').Value = '  +0.51%  '

# Row 46
$ws.Range('B46').Value = 'InjectiveProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D46').Value = "'28.10"
$ws.Range('E46').Value = '  -1.73%  '

# Row 47
$ws.Range('B47').Value = 'dogwifhat'
$ws.Range('C47').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D47').Value = "'2.68"
$ws.Range('E47').Value = '  -0.33%  '

# Row 48
$ws.Range('E48').Value = '  -2.53%  '

# Row 49
$ws.Range('D49').Value = "'0.000263"
$ws.Range('E49').Value = '  -5.12%  '

# Row 50
$ws.Range('D50').Value = "'7.74"
$ws.Range('E50').Value = '  -1.30%  '

# Row 51
$ws.Range('D51').Value = "'1.23"
$ws.Range('E51').Value = '  -3.92%  '
